# Update column F ("dSF") values on the active sheet per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -3
    4  = 3
    5  = -1
    7  = -5
    8  = 7
    9  = 3
    10 = 4
    11 = -4
    13 = -2
    14 = -4
    15 = -3
    16 = 3
    17 = -3
    18 = 5
    22 = -2
    23 = -2
    24 = -3
    25 = 2
    27 = 3
    28 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
